# Update (Removed Auto Arima)
# Updates forecast values on "Forecast Comparison" sheet (rows 2-17, columns C-G)
# and the derived totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New values for columns C (Prophet Forecast), D (Amazon Mean Forecast),
# E (Amazon P70 Forecast), F (Amazon P80 Forecast), G (Amazon P90 Forecast)
# for rows 2..17 (weeks W01..W16)
$data = @(
    @{ Row=2;  C=44; D=22; E=27; F=33; G=42 },
    @{ Row=3;  C=36; D=13; E=16; F=20; G=28 },
    @{ Row=4;  C=27; D=10; E=13; F=16; G=23 },
    @{ Row=5;  C=21; D=12; E=14; F=19; G=26 },
    @{ Row=6;  C=18; D=11; E=14; F=19; G=27 },
    @{ Row=7;  C=18; D=11; E=13; F=18; G=26 },
    @{ Row=8;  C=17; D=11; E=14; F=19; G=27 },
    @{ Row=9;  C=16; D=11; E=14; F=19; G=27 },
    @{ Row=10; C=16; D=11; E=14; F=18; G=27 },
    @{ Row=11; C=17; D=11; E=14; F=18; G=27 },
    @{ Row=12; C=19; D=11; E=14; F=19; G=27 },
    @{ Row=13; C=20; D=12; E=14; F=20; G=29 },
    @{ Row=14; C=18; D=12; E=14; F=19; G=29 },
    @{ Row=15; C=15; D=11; E=13; F=19; G=28 },
    @{ Row=16; C=12; D=11; E=13; F=18; G=28 },
    @{ Row=17; C=9;  D=10; E=13; F=18; G=26 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $wsForecast.Cells.Item($r, 3).Value = $entry.C   # Column C
    $wsForecast.Cells.Item($r, 4).Value = $entry.D   # Column D
    $wsForecast.Cells.Item($r, 5).Value = $entry.E   # Column E
    $wsForecast.Cells.Item($r, 6).Value = $entry.F   # Column F
    $wsForecast.Cells.Item($r, 7).Value = $entry.G   # Column G
}

# Update Summary sheet totals (stored as text in column B, same as before the
# edit). A leading apostrophe forces text storage instead of Excel
# auto-converting the numeric-looking string to a number; ClearFormats()
# afterwards drops the "quote prefix" cell style so no stray style is left
# behind on the cell.
$wsSummary.Range("B9").Value  = "'323"
$wsSummary.Range("B10").Value = "'197"
$wsSummary.Range("B11").Value = "'128"
$wsSummary.Range("B12").Value = "'44"
$wsSummary.Range("B14").Value = "'9"
$wsSummary.Range("B9:B14").ClearFormats()
